$d = $word.ActiveDocument

$pairs = @(
    @("78÷6=13, 0", "39÷9=4, 3"),
    @("20÷3=6, 2", "29÷6=4, 5"),
    @("20÷7=2, 6", "59÷4=14, 3"),
    @("50÷3=16, 2", "37÷3=12, 1"),
    @("10÷7=1, 3", "97÷4=24, 1"),
    @("38÷2=19, 0", "58÷4=14, 2"),
    @("40÷5=8, 0", "50÷5=10, 0"),
    @("98÷9=10, 8", "54÷5=10, 4"),
    @("76÷9=8, 4", "73÷6=12, 1"),
    @("27÷6=4, 3", "16÷5=3, 1"),
    @("71÷5=14, 1", "92÷6=15, 2"),
    @("73÷3=24, 1", "16÷3=5, 1"),
    @("78÷9=8, 6", "27÷7=3, 6"),
    @("62÷8=7, 6", "36÷4=9, 0"),
    @("24÷4=6, 0", "62÷9=6, 8"),
    @("10÷5=2, 0", "48÷4=12, 0"),
    @("51÷8=6, 3", "62÷8=7, 6"),
    @("72÷8=9, 0", "23÷3=7, 2"),
    @("44÷6=7, 2", "73÷7=10, 3"),
    @("65÷2=32, 1", "83÷9=9, 2"),
    @("35÷7=5, 0", "36÷9=4, 0"),
    @("95÷4=23, 3", "70÷2=35, 0"),
    @("29÷3=9, 2", "90÷2=45, 0"),
    @("10÷3=3, 1", "37÷7=5, 2"),
    @("75÷7=10, 5", "93÷2=46, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
